# Update "想去人数" (interest count) figures in F column on both the
# "展览" and "全部类型" sheets, as captured by the upstream data refresh
# (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row : new value } for column F updates.
$updates = @{
    "展览" = @{
        2  = 3078
        3  = 490
        9  = 1058
        10 = 14919
        14 = 5932
        18 = 89
        23 = 821
        26 = 10754
    }
    "全部类型" = @{
        3  = 3078
        4  = 490
        10 = 1058
        11 = 14919
        15 = 5932
        19 = 89
        24 = 821
        28 = 10754
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}

$wb.Save()
